$d = $word.ActiveDocument

$changes = @{
    5 = "[[PERSON_3]] – „pro [[PERSON_3]]“, „s [[PERSON_3]]“"
    6 = "[[PERSON_4]] – „s [[PERSON_4]]“, „o [[PERSON_4]]“"
    7 = "[[PERSON_5]] – „u [[PERSON_5]]“, „k [[PERSON_5]]“"
    8 = "[[PERSON_6]] – „o [[PERSON_6]]“, „se [[PERSON_6]]“"
    9 = "[[PERSON_7]] – „k [[PERSON_7]]“, „u [[PERSON_7]]“"
    10 = "[[PERSON_8]] – „s [[PERSON_8]]“, „o [[PERSON_8]]“"
    11 = "[[PERSON_9]] – „u [[PERSON_9]]“, „s [[PERSON_9]]“"
    12 = "[[PERSON_10]] – „s [[PERSON_10]]“, „k [[PERSON_10]]“"
    13 = "[[PERSON_11]] – „s [[PERSON_11]]“, „o [[PERSON_11]]“"
    14 = "[[PERSON_12]] – „ke [[PERSON_12]]“, „o [[PERSON_12]]“"
    15 = "[[PERSON_13]] – „o [[PERSON_13]]“, „s [[PERSON_13]]“"
    16 = "[[PERSON_14]] – „u [[PERSON_14]]“, „s [[PERSON_14]]“"
    17 = "[[PERSON_15]] – „ke [[PERSON_15]]“, „o [[PERSON_15]]“"
    18 = "[[PERSON_16]] – „s [[PERSON_16]]“, „o [[PERSON_16]]“"
    19 = "[[PERSON_17]] – „s [[PERSON_17]]“, „o [[PERSON_17]]“"
    20 = "[[PERSON_18]] – „k [[PERSON_18]]“, „od [[PERSON_18]]“"
    21 = "[[PERSON_19]] – „o [[PERSON_19]]“, „s [[PERSON_19]]“"
    22 = "[[PERSON_20]] – „o [[PERSON_20]]“, „se [[PERSON_20]]“"
    23 = "[[PERSON_21]] – „s [[PERSON_21]]“, „u [[PERSON_21]]“"
    24 = "[[PERSON_22]] – „o [[PERSON_22]]“, „s [[PERSON_22]]“"
    25 = "[[PERSON_23]] – „k [[PERSON_23]]“, „o [[PERSON_23]]“"
    26 = "[[PERSON_24]] – „se [[PERSON_24]]“, „o Soně Mikulkové“"
    27 = "[[PERSON_25]] – „o [[PERSON_25]]“, „s [[PERSON_25]]“"
    30 = "[[PERSON_26]] – „s [[PERSON_26]]“, „o [[PERSON_26]]“"
    31 = "[[PERSON_27]] – „k [[PERSON_27]]“, „s [[PERSON_27]]“"
    32 = "[[PERSON_28]] – „s [[PERSON_28]]“, „o [[PERSON_28]]“"
    33 = "[[PERSON_29]] – „od [[PERSON_29]]“, „s [[PERSON_29]]“"
    34 = "[[PERSON_30]] – „k [[PERSON_31]]“, „o [[PERSON_31]]“"
    35 = "[[PERSON_32]] – „o [[PERSON_32]]“, „s [[PERSON_32]]“"
    36 = "[[PERSON_33]] – „s [[PERSON_33]]“, „o [[PERSON_33]]“"
    37 = "[[PERSON_34]] – „s [[PERSON_34]]“, „o [[PERSON_34]]“"
    38 = "[[PERSON_35]] – „k [[PERSON_35]]“, „s [[PERSON_35]]“"
    39 = "[[PERSON_36]] – „pro [[PERSON_36]]“, „o [[PERSON_36]]“"
    40 = "[[PERSON_37]] – „k [[PERSON_37]]“, „o [[PERSON_37]]“"
    41 = "[[PERSON_38]] – „o [[PERSON_38]]“, „s [[PERSON_38]]“"
    42 = "[[PERSON_39]] – „s [[PERSON_39]]“, „o [[PERSON_39]]“"
    43 = "[[PERSON_40]] – „s [[PERSON_40]]“, „o [[PERSON_40]]“"
    44 = "[[PERSON_41]] – „u [[PERSON_41]]“, „o [[PERSON_41]]“"
    45 = "[[PERSON_42]] – „se [[PERSON_42]]“, „o [[PERSON_42]]“"
    46 = "[[PERSON_43]] – „o [[PERSON_43]]“, „s [[PERSON_43]]“"
    47 = "[[PERSON_44]] – „k [[PERSON_44]]“, „o [[PERSON_44]]“"
    48 = "[[PERSON_45]] – „o [[PERSON_45]]“, „s [[PERSON_45]]“"
    49 = "[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_46]]“"
    50 = "[[PERSON_47]] – „s [[PERSON_47]]“, „o [[PERSON_47]]“"
    51 = "[[PERSON_48]] – „o [[PERSON_48]]“, „s [[PERSON_48]]“"
    52 = "[[PERSON_49]] – „s [[PERSON_49]]“, „o [[PERSON_49]]“"
    53 = "[[PERSON_50]] – „o [[PERSON_50]]“, „s [[PERSON_50]]“"
    54 = "[[PERSON_51]] – „s [[PERSON_51]]“, „o [[PERSON_51]]“"
}

foreach ($idx in $changes.Keys) {
    $p = $d.Paragraphs.Item([int]$idx)
    $p.Range.Text = $changes[$idx]
}

Write-Output ("Updated " + $changes.Count + " paragraphs")
